# Generate Report for Handback
# Update the "Latest Handback DateTime" (column K, row 2) on the zh-cn and
# de-de sheets to reflect newly generated handback timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-10-18 05:00:01"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-10-18 05:00:27"
